$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,('Designator', 'Description', 'Comment', 'Manufacturer', 'Manufacturer Part Number', 'Digi-Key Part Number', 'PackageReference', 'Quantity', 'GENERIC_OK')
    ,('C1, C2, C6, C7, C11', '1µF ±10% 10V X5R Ceramic Capacitor 0402', '1µF-10V', 'Yageo', 'CC0402KRX5R6BB105', '311-1439-1-ND', '0402', 5, '')
    ,('C3, C8, C9, C10', '4.7µF ±10% 6.3V X5R Ceramic Capacitor 0402', '4.7µF', 'TDK Corporation', 'C1005X5R0J475K050BC', '445-5947-1-ND', '0402', 4, '')
    ,('C4, C5', '16µF ±10% 16V X5R Ceramic Capacitor 0603', '10µF-16V', 'Murata Electronics North America', 'GRT188R61C106KE13D', '490-12317-1-ND', '0603', 2, '')
    ,('C12', '2.2nF ±10% 50V X5R Ceramic Capacitor 0402', '2.2nF', 'Murata Electronics North America', 'GRM155R61H222KA01D', '490-8257-1-ND', '0402', 1, '')
    ,('C13, C14', 'CAP CER 50PF 50V NP0 0603', '50pF', 'AVX Corporation', '06035A500JAT2A', '478-10324-1-ND', '0603', 2, '')
    ,('C15, C16', '0.1µF ±10% 10V X5R Ceramic Capacitor 0402', '0.1µF-10V', 'Yageo', 'CC0402KRX5R6BB104', '311-1336-1-ND', '0402', 2, '')
    ,('C17', '1.2nF ±10% 10V U2J Ceramic Capacitor 0402', '1.2nF-10V', 'Murata Electronics North America', 'GRM1557U1A122JA01D', '490-6272-1-ND', '0402', 1, '')
    ,('CON1', 'USB Micro B', 'UJ2-MIBH-4-MSMT-TR', 'CUI Inc.', 'UJ2-MIBH-4-MSMT-TR', '102-4008-1-ND', '', 1, '')
    ,('CON2', '1.20MM WTB HEADER VERT 2POS', '0781710002', 'Molex, LLC', '0781710002', 'WM4459CT-ND', '', 1, '')
    ,('D1', '4.7µF ±10% 6.3V X5R Ceramic Capacitor 0402', 'LXES15AAA1-133', 'Murata Electronics North America', 'LXES15AAA1-133', '490-13209-1-ND', '0402', 1, '')
    ,('D2', '0603 INFRARED DIODE', 'APT1608F3C', 'Kingbright', 'APT1608F3C', '754-1433-1-ND', '0603', 1, '')
    ,('J2', 'MICRO UNIVERSAL CONTACT Z 3.5MM', '120220-0314', 'ITT Cannon, LLC', '120220-0314', '1003-2383-1-ND', '', 1, '')
    ,('L1', 'FIXED IND 15UH 350MA 0.6 OHM SMD', '15µH', 'TDK Corporation', 'MLZ2012N150LT000', '445-6763-1-ND', '0805', 1, '')
    ,('L2', 'FIXED IND 2.2UH 600MA 0.2OHM SMD', '2.2µH', 'TDK Corporation', 'MLZ2012M2R2HTD25', '445-17077-1-ND', '0805', 1, '')
    ,('L3', '0805 Ferrite Bead 330Ohm@100MHz', 'BLM21PG331SN1D', 'Murata Electronics North America', 'BLM21PG331SN1D', '490-5988-1-ND', '0805 (2012 Metric)', 1, '')
    ,('MOD1', 'BLE PSoC Module 4.1', 'CYBLE-214009-00', 'Cypress Semiconductor Corp', 'CYBLE-214009-00', '428-3809-ND', '', 1, '')
    ,('Q1', '0805 PHOTOTRANSISTOR', 'APT2012P3BT', 'Kingbright', 'APT2012P3BT', '754-1436-1-ND', '0805', 1, '')
    ,('R1', 'RES SMD 820kOHM 1% 1/16W 0402', '820K', 'Yageo', 'RC0402FR-07820KL', 'YAG3240CT-ND', '0402', 1, 'YES')
    ,('R2', 'RES SMD 82kOHM 1% 1/16W 0402', '82K', 'Yageo', 'RC0402FR-0782KL', '311-82.0KLRCT-ND', '0402', 1, 'YES')
    ,('R3, R4', 'RES SMD 200kOHM 1% 1/16W 0402', '200K', 'Yageo', 'AC0402FR-07200KL', 'YAG3452CT-ND', '0402', 2, 'YES')
    ,('R5, R13, R14, R17', 'RES SMD 10kOHM 1% 1/16W 0402', '10K', 'Yageo', 'RC0402FR-0710KL', '311-10.0KLRCT-ND', '0402', 4, 'YES')
    ,('R6', 'RES SMD 5.6kOHM 5% 1/16W 0402', '5.6K', 'Yageo', 'RC0402JR-075K6L', '311-5.6KJRCT-ND', '0402', 1, 'YES')
    ,('R7, R11, R12', 'RES SMD 41.2kOHM 1% 1/16W 0402', '41.2K', 'Yageo', 'RC0402FR-0741K2L', 'YAG3155CT-ND', '0402', 3, 'YES')
    ,('R8, R9, R18, R19', 'Do Not Place', 'Do Not Place', '', '', '', '', 4, '')
    ,('R10', 'RES SMD 390kOHM 1% 1/16W 0402', '390K', 'Yageo', 'RC0402FR-07390KL', 'YAG3143CT-ND', '0402', 1, 'YES')
    ,('R15', 'RES SMD 100OHM 1% 1/16W 0402', '100', 'Yageo', 'RC0402FR-07100RL', '311-100LRCT-ND', '0402', 1, 'YES')
    ,('R16', 'RES SMD 0.0OHM JUMPER 1/16W 0402', '0', 'Yageo', 'RC0402JR-070RL', '311-0.0JRCT-ND', '0402', 1, 'YES')
    ,('SW1, SW2', 'SPST Tactile Switch', 'KXT311LHS', 'C&K', 'KXT 311 LHS', 'CKN10777CT-ND', '', 2, '')
    ,('U1', 'Boost converter 28V', 'TPS61046YFFR', 'Texas Instruments', 'TPS61046YFFR', '296-44337-1-ND', '6-UFBGA, DSBGA', 1, '')
    ,('U2', 'Buck Converter', 'TPS627431YFPR', 'Texas Instruments', 'TPS627431YFPR', '296-44499-2-ND', '8-XFBGA, DSBGA', 1, '')
    ,('U3', 'Li-Ion linear battery charger', 'STBC03', 'STMicroelectronics', 'STBC03JR', '497-17241-1-ND', '30-UFBGA, FCBGA', 1, '')
    ,('U5', 'IC EEPROM 64KBIT 400KHZ 6WLCSP', 'LE2464CXATBG', 'ON Semiconductor', 'LE2464CXATBG', 'LE2464CXATBGOSDKR-ND', '6-XFBGA, WLCSP', 1, '')
    ,('U6', 'IC I/O EXPANDER I2C 12XQFN', 'PCA9571GUX', 'NXP USA Inc.', 'PCA9571GUX', '568-12400-1-ND', '12-XFQFN', 1, '')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 1
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $row[$j]
    }
}